# "Add a walk with the dog"
# Bump February's extra walking distance (Sheet1!G3) by 2.3, which ripples
# through the running-total formula in F3 and the "Actual" series cache on
# the Chart1 sheet.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
$sheet1.Range("G3").Value = 66.8
$sheet1.Range("G4").Select()

$chartSheet = $wb.Worksheets.Item("Chart1")
$chartSheet.Range("A23").Select()
